# "Delete Post, User Wall & Friends Posts , Comments functionality added,
#  Bugs Fixed, News Feed Page implemented"
#
# The self-evaluation scoring column (C) was filled in with the actual
# scores achieved for each criterion (it was previously blank / not
# scored yet). One cell (C11, "Web Design") had held a placeholder single
# space character as text; it now gets its real numeric score instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Self-Evaluation-Protocol")

# Row -> score achieved in column C.
$scores = [ordered]@{
    8  = 9    # Days Commit in GitHub
    9  = 15   # Numbers of Commits in GitHub
    11 = 7    # Web Design (was a literal " " text placeholder)
    12 = 30   # AngularJS Project Structure
    13 = 5    # Login Screen
    14 = 10   # User Register Screen
    15 = 5    # User Home Screen
    16 = 10   # -->Search by Username
    17 = 3    # -->Display Pending Requests
    18 = 7    # -->Display Requests Details
    19 = 10   # -->Display Top Friends with Images
    20 = 5    # -->Display Posts by Friends (Feed)
    21 = 10   # -->Display Post Data
    22 = 10   # -->Display Dropdown for Commenting
    23 = 5    # -->Display Like/Unlike Buttons and Logic
    25 = 10   # User Wall
    26 = 10   # -->Post Box with Submit Button
    27 = 10   # Friends
    29 = 5    # Delete Post
    30 = 10   # Edit User Profile
    31 = 5    # Change User Password
    32 = 5    # Logout
    33 = 5    # Guest Authorization Checks
    34 = 10   # User Authorization Checks
    42 = 20   # User Authorization Checks For Comments
}

foreach ($row in $scores.Keys) {
    $ws.Cells.Item($row, 3).Value = $scores[$row]
}

# Recalculate the Total Score formula (C44) and move the cursor/selection
# to where the author left off (row 2 scrolled into view, C41 selected).
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("C41").Select()

$excel.Calculate()
